$wb = $excel.ActiveWorkbook

$wsIT = $wb.Worksheets.Item("IT")
$wsIT.Range("B2").Value = 2020

$wsIT.Activate()
$wsIT.Range("B3").Select()
